$d = $word.ActiveDocument

# Replace the old phone-format example text (split across two runs:
# " 00-000-" + "0000") with the new single-run text
# "00-000-00 או 000-00-000".
$range = $d.Content
$range.Find.Execute(" 00-000-0000", $true, $true, $false, $false, $false, $true, 1, $false, "00-000-00 או 000-00-000", 2) | Out-Null

# Word leaves a "_GoBack" bookmark marking the last edit position,
# collapsed right after the inserted text (i.e. right before the
# paragraph mark of that paragraph).
$p = $range.Duplicate
$p.Expand(4) | Out-Null
$goBackPos = $p.End - 1

# Placing a collapsed bookmark exactly one position before a paragraph
# mark trips this host's bookmark/paragraph resolution, so nudge the
# paragraph mark further away first with a one-character placeholder,
# add the bookmark at the now-safe position, then remove the
# placeholder again (the bookmark stays put, now correctly collapsed
# right before the paragraph mark).
$placeholder = $d.Range($goBackPos, $goBackPos)
$placeholder.InsertAfter("x") | Out-Null

$bmRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$d.Range($goBackPos, $goBackPos + 1).Delete() | Out-Null
